$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# "Enterprises density (per 1000 people)" row (row 11): Micro, SMEs, MSMEs
# These are stored as text values (not numbers), so a leading apostrophe
# forces Excel to keep them as text instead of converting to numeric.
$ws.Range("B11").Value = "'28.59"
$ws.Range("C11").Value = "'4.39"
$ws.Range("D11").Value = "'32.98"

# "Enterprises (% of total)" row (row 12): SMEs, MSMEs (Micro column stays 86.2)
$ws.Range("C12").Value = "'13.22"
$ws.Range("D12").Value = "'99.43"
